# Update countries & provincias Spain
# Applies the covid-data refresh captured in the commit diff:
#  - timestamp label bumped from 06:32 to 07:49
#  - several countries' stats refreshed, which causes a few adjacent
#    countries to swap rank/row because one overtook the other in total cases
#  - plain numeric refreshes for a handful of other countries that did not
#    change rank

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header timestamp (A1)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 11 de Agosto de 2020 a las 07:49"

# ---------------------------------------------------------------------
# 2) Rank swaps: the country whose numbers grew moves up a row, pushing
#    its neighbour down one row with its old (unchanged) numbers.
# ---------------------------------------------------------------------

# India .. (unchanged rank) simple refresh is handled below in section 3.

# Armenia (row55) / Kirguistan (row56) -> Kirguistan overtakes Armenia
$ws.Range("A55").Value = "Kirguistan"
$ws.Range("B55").Value = 40455
$ws.Range("C55").Value = 278
$ws.Range("D55").Value = 32734
$ws.Range("E55").Value = 6243
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 4
$ws.Range("H55").Value = 1478

$ws.Range("A56").Value = "Armenia"
$ws.Range("B56").Value = 40433
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 32616
$ws.Range("E56").Value = 7021
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 796

# Polinesia Francesa (row191) / Butan (row192) -> Butan overtakes Polinesia Francesa
$ws.Range("A191").Value = "Butan"
$ws.Range("B191").Value = 113
$ws.Range("C191").Value = 3
$ws.Range("D191").Value = 97
$ws.Range("E191").Value = 16
$ws.Range("F191").Value = 0
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 0

$ws.Range("A192").Value = "Polinesia Francesa"
$ws.Range("B192").Value = 112
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 62
$ws.Range("E192").Value = 50
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 0

# Santa Lucia (row202) / Timor Oriental (row203) swap places (numbers tied,
# so no visible numeric change, just the country labels trade rows)
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("B202").Value = 25
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 24
$ws.Range("E202").Value = 1
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0

$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("B203").Value = 25
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 24
$ws.Range("E203").Value = 1
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

# Islas Malvinas (row213) / Montserrat (row214) -> Montserrat overtakes Islas Malvinas
$ws.Range("A213").Value = "Montserrat"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 12
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

# ---------------------------------------------------------------------
# 3) Plain numeric refreshes (rank / country unchanged)
# ---------------------------------------------------------------------

# India (row6)
$ws.Range("B6").Value = 2269052
$ws.Range("C6").Value = 1899
$ws.Range("D6").Value = 1583489
$ws.Range("E6").Value = 640202
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 8
$ws.Range("H6").Value = 45361

# Pakistan (row17)
$ws.Range("B17").Value = 285191
$ws.Range("C17").Value = 531
$ws.Range("D17").Value = 261246
$ws.Range("E17").Value = 17833
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 6112

# Israel (row33)
$ws.Range("B33").Value = 85222
$ws.Range("C33").Value = 500
$ws.Range("D33").Value = 59999
$ws.Range("E33").Value = 24610
$ws.Range("F33").Value = 0

# Uzbekistan (row62)
$ws.Range("B62").Value = 31545
$ws.Range("C62").Value = 241
$ws.Range("D62").Value = 22992
$ws.Range("E62").Value = 8351
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 202

# Australia (row72)
$ws.Range("B72").Value = 21713
$ws.Range("C72").Value = 316
$ws.Range("D72").Value = 12144
$ws.Range("E72").Value = 9238
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 18
$ws.Range("H72").Value = 331

# El Salvador (row73)
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 9762
$ws.Range("E73").Value = 10540
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 7
$ws.Range("H73").Value = 570

# Tailandia (row115)
$ws.Range("D115").Value = 3163
$ws.Range("E115").Value = 130

# Camboya (row179)
$ws.Range("B179").Value = 266
$ws.Range("C179").Value = 15
$ws.Range("D179").Value = 220
$ws.Range("E179").Value = 46
